# Regenerate save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 6
    3  = 6
    4  = 2
    5  = 5
    6  = 6
    7  = 1
    8  = 9
    9  = 9
    10 = 6
    11 = 4
    12 = 4
    13 = 1
    14 = 3
    15 = 3
    16 = 12
    17 = 5
    18 = 5
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
